# Update "想去人数" (attendance interest count) values that changed
# between the previous data scrape and the new one (456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14193
$ws1.Range("F4").Value = 686
$ws1.Range("F6").Value = 549
$ws1.Range("F7").Value = 1485
$ws1.Range("F8").Value = 139

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14193
$ws4.Range("F4").Value = 686
$ws4.Range("F8").Value = 549
$ws4.Range("F9").Value = 1485
$ws4.Range("F11").Value = 139
